# Update export pin format
#
# The old layout had a 2-row, green "banner" title (B2:F3 merged) followed by a
# blank spacer row, then the "Exam Title:" / "Start  Date:" label rows and a
# trailing blank row before the data-table header row.
#
# The new layout keeps the green banner title but shrinks it to a single
# (taller) row, adds a plain blank spacer row under it, shifts the label rows
# up to follow immediately, and adds one extra bold/left-aligned blank row
# just above the data-table header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second half of the old title banner (row 3). Excel shifts every
# row below it up by one, which automatically:
#   - turns the old "B2:F3" banner merge into "B2:F2"
#   - turns the old "B5:C5" ("Exam Title:") merge into "B4:C4"
#   - moves "Start  Date:" from row 6 to row 5, and the blank row from 7 to 6
$ws.Rows(3).Delete()

# The single-row banner now needs to be taller than the original first half.
$ws.Rows(2).RowHeight = 29.25

# Insert a new blank row just above the data-table header (now row 7, since
# everything shifted up by one already) and give it the small spacer height.
$ws.Rows(7).Insert()
$ws.Rows(7).RowHeight = 12.75

# Style the new blank cell bold + left aligned (matching the sheet's other
# label cells), leaving it otherwise empty.
$ws.Range("B7").Font.Bold = $true
$ws.Range("B7").HorizontalAlignment = -4131

# Put the active selection on E4 (the "Exam Title" value input cell).
$null = $ws.Range("E4").Select()
